$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C) for rows 2-9 from 45170 to 45174
# (date serial number change, e.g. 2023-09-01 -> 2023-09-05)
foreach ($row in 2..9) {
    $ws.Cells.Item($row, 3).Value = 45174
}
